$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (avoid numeric auto-coercion of values like "1.013")
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.935.72'
$ws.Range("E2").Value = '  -2.09%  '

$ws.Range("D3").Value = '1.787.41'
$ws.Range("E3").Value = '  -2.42%  '

$ws.Range("D4").Value = '1.013'
$ws.Range("E4").Value = '  +1.00%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '311.73'
$ws.Range("E5").Value = '  -0.86%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.91%  '

$ws.Range("D7").Value = '0.4230'
$ws.Range("E7").Value = '  -1.57%  '

$ws.Range("D8").Value = '0.3603'
$ws.Range("E8").Value = '  -1.57%  '

$ws.Range("D9").Value = '0.07163'
$ws.Range("E9").Value = '  -1.62%  '

$ws.Range("D10").Value = '0.8413'
$ws.Range("E10").Value = '  -3.43%  '

$ws.Range("D11").Value = '20.24'
$ws.Range("E11").Value = '  -2.11%  '

$ws.Range("D12").Value = '1.896.57'
$ws.Range("E12").Value = '  +4.03%  '

$ws.Range("D13").Value = '5.260'

$ws.Range("D14").Value = '6.352'
$ws.Range("E14").Value = '  -2.85%  '

$ws.Range("D15").Value = '0.06845'
$ws.Range("E15").Value = '  -1.39%  '

$ws.Range("E16").Value = '  +0.89%  '

$ws.Range("D17").Value = '79.73'
$ws.Range("E17").Value = '  -0.99%  '

$ws.Range("D18").Value = '0.000008696'
$ws.Range("E18").Value = '  -2.68%  '

$ws.Range("D19").Value = '1.010'

$ws.Range("D20").Value = '14.92'
$ws.Range("E20").Value = '  -3.47%  '

$ws.Range("D21").Value = '27.121.99'
$ws.Range("E21").Value = '  -1.30%  '

$ws.Range("D22").Value = '5.044'
$ws.Range("E22").Value = '  -2.41%  '

$ws.Range("E23").Value = '  +1.86%  '

$ws.Range("D24").Value = '2.062.70'
$ws.Range("E24").Value = '  +1.79%  '

$ws.Range("D25").Value = '1.950'
$ws.Range("E25").Value = '  -1.49%  '

$ws.Range("D26").Value = '153.30'
$ws.Range("E26").Value = '  -0.92%  '

$ws.Range("D27").Value = '18.19'
$ws.Range("E27").Value = '  -3.62%  '

$ws.Range("D28").Value = '114.73'
$ws.Range("E28").Value = '  +0.22%  '

$ws.Range("D29").Value = '5.007'
$ws.Range("E29").Value = '  -3.26%  '

$ws.Range("E30").Value = '  -12.12%  '

$ws.Range("D31").Value = '0.08930'
$ws.Range("E31").Value = '  +0.27%  '

$ws.Range("D32").Value = '0.7224'
$ws.Range("E32").Value = '  -5.15%  '

$ws.Range("D33").Value = '2.849'
$ws.Range("E33").Value = '  -3.95%  '

$ws.Range("D34").Value = '4.316'
$ws.Range("E34").Value = '  -5.22%  '

$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").Value = '1.010'
$ws.Range("E35").Value = '  +0.75%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.082'
$ws.Range("E36").Value = '  -5.40%  '

$ws.Range("D37").Value = '1.082'
$ws.Range("E37").Value = '  -1.31%  '

$ws.Range("D38").Value = '0.01894'
$ws.Range("E38").Value = '  -2.37%  '

$ws.Range("D39").Value = '0.05082'
$ws.Range("E39").Value = '  -4.54%  '

$ws.Range("D40").Value = '0.4945'
$ws.Range("E40").Value = '  -2.98%  '

$ws.Range("D41").Value = '0.1614'
$ws.Range("E41").Value = '  -3.28%  '

$ws.Range("D42").Value = '2.509'
$ws.Range("E42").Value = '  -10.60%  '

$ws.Range("D43").Value = '5.978'
$ws.Range("E43").Value = '  -9.80%  '

$ws.Range("D44").Value = '7.943'
$ws.Range("E44").Value = '  -5.92%  '

$ws.Range("D45").Value = '1.012'
$ws.Range("E45").Value = '  +1.08%  '

$ws.Range("D46").Value = '104.47'
$ws.Range("E46").Value = '  -1.67%  '

$ws.Range("D47").Value = '10.06'
$ws.Range("E47").Value = '  -4.42%  '

$ws.Range("D48").Value = '0.06286'
$ws.Range("E48").Value = '  -3.44%  '

$ws.Range("D49").Value = '0.4478'
$ws.Range("E49").Value = '  -4.68%  '

$ws.Range("D50").Value = '1.578'
$ws.Range("E50").Value = '  -2.77%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '62.33'
$ws.Range("E51").Value = '  -3.06%  '

# Clean up the temporary text formatting so cells fall back to the default style
$ws.Range("D2:D51").Style = "Normal"

